# Update match #7 (row 8) with its final result: Am. Ma. Velho 5 x 0 Tira Fama
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("partidas")

$ws.Range("E8").Value = "5x0"
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = "Finalizado"

# Move the active selection like the author left it after editing
$ws.Range("L8").Select()
